# Generate Report for Handoff
# Updates status/datetime/error-detail info for the b8174788-... file
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c2c0219b9058da9d087bad0c42b6fcbb7c174c0/e2e/b8174788-9a90-4227-8136-0f93962b431e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09571c54074ba7cfc7e46f0826786c3498bced26/e2e/b8174788-9a90-4227-8136-0f93962b431e.md."

# Overview sheet - row 3 corresponds to b8174788-...md
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 06:47:42"

# zh-cn sheet - row 3 corresponds to b8174788-...md
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-06 06:47:31"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet - row 3 corresponds to b8174788-...md
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-06 06:47:42"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
